$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Данные"

# Update data values (B5:B11 divided by 10)
$ws.Range("B5").Value = 14000
$ws.Range("B6").Value = 13600
$ws.Range("B7").Value = 84000
$ws.Range("B8").Value = 42000
$ws.Range("B9").Value = 11800
$ws.Range("B10").Value = 10400
$ws.Range("B11").Value = 11200

# Fix F31's formatting so it matches the rest of column F (it had a stray
# duplicate style); copy formatting from a neighboring F-column cell.
$ws.Range("F30").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Change selection
$ws.Range("D13").Select()
